# Add "Source" footer block to the worksheet, reflecting the
# "update data downloads for reform with sources" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 20-22: a "Source:" label (styled like the other headers)
# followed by the data source name and URL.
$ws.Range("A20").Value = "Source:"
$ws.Range("A21").Value = "National Corrections Reporting Program"
$ws.Range("A22").Value = "https://www.bjs.gov/index.cfm?ty=dcdetail&iid=268"

# Give A20 the same header styling used by A1/A3 ("Year" style).
$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial(-4122)

# Reflect the user's resulting selection/scroll position after the edit.
$ws.Range("A20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 14
